$wb = $excel.ActiveWorkbook

# --- 1. Rename "Requested quantity" header on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet right after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Forecast data rows
$forecastData = @(
  @(45564.99999999999, 389, -277.9026598592412, 1001.56888963263),
  @(45578.99999999999, 402, -233.3172966593288, 1044.412280409216),
  @(45585.99999999999, 408, -233.8048571030846, 1089.06677877703),
  @(45592.99999999999, 414, -260.1288001162836, 1028.330496663663),
  @(45599.99999999999, 421, -197.8365305642837, 1048.630268425087),
  @(45606.99999999999, 427, -219.5501954072548, 1081.307291174481),
  @(45620.99999999999, 440, -197.1934662729422, 1073.335168097889),
  @(45627.99999999999, 446, -204.6801817754492, 1039.171772024048),
  @(45634.99999999999, 453, -188.9747188899827, 1090.869377239449),
  @(45641.99999999999, 459, -182.0605373835199, 1109.057406904045),
  @(45648.99999999999, 465, -138.7903252124046, 1127.633980533924),
  @(45655.99999999999, 472, -170.805277107265, 1056.035828010919),
  @(45662.99999999999, 478, -189.6709322259946, 1102.572895242658),
  @(45669.99999999999, 485, -153.4426154225741, 1115.171621552917),
  @(45676.99999999999, 491, -125.0262868403675, 1154.210550117814)
)

$row = 2
foreach ($entry in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $entry[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $entry[1]
    $wsForecast.Cells.Item($row, 3).Value = $entry[2]
    $wsForecast.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
